$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.060.04'
$ws.Range('E2').Value = '  +0.05%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.367.32'
$ws.Range('E3').Value = '  +1.80%  '

$ws.Range('E4').Value = '  -0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.37'
$ws.Range('E5').Value = '  +0.52%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.25'
$ws.Range('E6').Value = '  +1.30%  '

$ws.Range('E7').Value = '  -0.35%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('E9').Value = '  -2.31%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.33'
$ws.Range('E10').Value = '  +1.03%  '

$ws.Range('E11').Value = '  +0.58%  '

$ws.Range('E12').Value = '  +1.49%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.58'
$ws.Range('E13').Value = '  -1.29%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.72'
$ws.Range('E14').Value = '  +0.85%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.736.98'
$ws.Range('E15').Value = '  +1.39%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.356.31'
$ws.Range('E16').Value = '  +1.64%  '

$ws.Range('E17').Value = '  +1.21%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.030.76'
$ws.Range('E18').Value = '  +0.16%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.99'
$ws.Range('E19').Value = '  -0.38%  '

$ws.Range('E20').Value = '  +2.25%  '

$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.10'
$ws.Range('E22').Value = '  +0.23%  '

$ws.Range('E23').Value = '  -0.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('E24').Value = '  -2.14%  '

$ws.Range('E25').Value = '  -6.63%  '

$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('E27').Value = '  -0.13%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.38'
$ws.Range('E28').Value = '  +7.65%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.33'
$ws.Range('E29').Value = '  +2.55%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.36'
$ws.Range('E30').Value = '  +2.86%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.30%  '

$ws.Range('E32').Value = '  +0.85%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.52'
$ws.Range('E33').Value = '  -0.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0729'
$ws.Range('E34').Value = '  +4.42%  '

$ws.Range('E35').Value = '  +6.75%  '

$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '127.95'
$ws.Range('E36').Value = '  -8.44%  '

$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('E37').Value = '  +0.97%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.33'
$ws.Range('E38').Value = '  -0.83%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.85'
$ws.Range('E39').Value = '  +3.74%  '

$ws.Range('E40').Value = '  -1.37%  '

$ws.Range('E41').Value = '  -0.57%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.80'
$ws.Range('E42').Value = '  -6.67%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.934.98'
$ws.Range('E43').Value = '  -0.21%  '

$ws.Range('E44').Value = '  +0.11%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.14'
$ws.Range('E45').Value = '  +3.86%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.30'
$ws.Range('E46').Value = '  -8.88%  '

$ws.Range('E47').Value = '  +0.23%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.601.74'
$ws.Range('E48').Value = '  +1.43%  '

$ws.Range('E49').Value = '  +2.77%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.75'
$ws.Range('E50').Value = '  -0.40%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.29'
$ws.Range('E51').Value = '  -2.60%  '
